$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115 - this pushes the existing rows 115..147 down to 116..148,
# carrying all their original values/formatting with them.
$ws.Rows.Item(115).Insert()

# Populate the newly-inserted row 115 with this week's record (same shape as the
# surrounding rows, new date / volume / price / origin for this observation).
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value = "Ñuble"
$ws.Cells.Item(115, 4).Value = 45093
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100108
$ws.Cells.Item(115, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(115, 9).Value = 100108002
$ws.Cells.Item(115, 10).Value = "Mango"
$ws.Cells.Item(115, 11).Value = "Sin especificar"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 70
$ws.Cells.Item(115, 14).Value = 9000
$ws.Cells.Item(115, 15).Value = 10000
$ws.Cells.Item(115, 16).Value = 9571
$ws.Cells.Item(115, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(115, 18).Value = "Perú"
$ws.Cells.Item(115, 19).Value = 2393
$ws.Cells.Item(115, 20).Value = 4
